$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = $ws.Range("C3").Value()
$ws.Range("E4").Value = $ws.Range("C4").Value()
$ws.Range("E5").Value = $ws.Range("C5").Value()
$ws.Range("E6").Value = $ws.Range("C6").Value()
$ws.Range("E7").Value = $ws.Range("C7").Value()

$ws.Range("F3").Select()
